# "test to get cell color" - add a few sample rows below the existing
# value/info header, each showing a different way a cell can carry color:
# a border, a plain RGB fill, a theme-based fill, and another RGB fill.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "int, border" label + number with a medium box border
$ws.Range("A2").Value = "int, border"
$ws.Range("B2").Value = 12
$ws.Range("B2").Borders.Weight = -4138   # xlMedium

# Row 3: "int, bgcolor" label + number with a yellow fill
$ws.Range("A3").Value = "int, bgcolor"
$ws.Range("B3").Value = 15
$ws.Range("B3").Interior.Color = 65535   # RGB(255,255,0) yellow

# Row 4: "int" label + number with a theme-colored fill + "theme" label
$ws.Range("A4").Value = "int"
$ws.Range("B4").Value = 17
# Priming Interior.Color first makes the pattern already "solid" so the
# following ThemeColor assignment doesn't leave a stray intermediate fill
# behind in the style table.
$ws.Range("B4").Interior.Color = 65535
$ws.Range("B4").Interior.ThemeColor = 10 # theme color index 9 (accent6)
$ws.Range("C4").Value = "theme"

# Row 5: number with a red fill
$ws.Range("B5").Value = 20
$ws.Range("B5").Interior.Color = 255     # RGB(255,0,0) red

$ws.Range("B5").Select()

